# Data Setup For Product Historical
# Reshape the flat single-row currency/product test data on "AddProductCategory1"
# into a proper 3-row table (one row per currency: USD / CAD / EUR).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddProductCategory1")
$ws.Activate() | Out-Null

# Row 1 gets a touch taller.
$ws.Rows.Item(1).RowHeight = 16.25

# J1/K1 used to hold "Canadian Dollar" / 11 - they now hold the numeric series
# that used to live in O1/S1.
$ws.Range("J1").Value = 20
$ws.Range("K1").Value = 30

# The rest of the old single-row spill (L1, M1, O1, P1, Q1, S1, T1, U1) is gone -
# its contents move onto rows 2 and 3 below.
$ws.Range("L1").ClearContents() | Out-Null
$ws.Range("M1").ClearContents() | Out-Null
$ws.Range("O1").ClearContents() | Out-Null
$ws.Range("P1").ClearContents() | Out-Null
$ws.Range("Q1").ClearContents() | Out-Null
$ws.Range("S1").ClearContents() | Out-Null
$ws.Range("T1").ClearContents() | Out-Null
$ws.Range("U1").ClearContents() | Out-Null

# R1 keeps its date formatting but the actual date value (43831) moves down to
# G2, so just clear the value and leave the formatted, empty cell behind.
$ws.Range("R1").ClearContents() | Out-Null

# Row 2: Canadian Dollar line.
$ws.Range("G2").Value = 43831
$ws.Range("G2").NumberFormat = "DD/MM/YY"
$ws.Range("H2").Value = "Canadian Dollar"
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 21
$ws.Range("K2").Value = 31

# Row 3: Euro line.
$ws.Range("H3").Value = "Euro"
$ws.Range("I3").Value = 12
$ws.Range("J3").Value = 22
$ws.Range("K3").Value = 32

# Leave the selection where the author left it.
$ws.Range("G6").Select() | Out-Null
